# Applies the cryptos.xlsx crypto-price/volume refresh described in the commit diff.
# Every cell written here (B/C/D/E for rows 2-51) is text data in the workbook (coin
# names, links, formatted price strings and padded percentage strings), so plain numeric
# looking values (e.g. "1.00", "6.20", "557.42") are assigned with a leading apostrophe
# to force Excel to keep them as text instead of silently coercing them to floating point
# numbers (which would corrupt values like "1.00" -> 1 or "557.42" -> 557.41999999999996).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.621.44'
$ws.Range('E2').Value = '  -1.54%  '
$ws.Range('D3').Value = '3.045.25'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'557.42"
$ws.Range('E5').Value = '  -0.35%  '
$ws.Range('D6').Value = "'141.87"
$ws.Range('E6').Value = '  -1.38%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '3.043.64'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').Value = "'0.518"
$ws.Range('E9').Value = '  +3.88%  '
$ws.Range('E10').Value = '  -10.33%  '
$ws.Range('D11').Value = "'0.153"
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = "'0.487"
$ws.Range('E12').Value = '  +5.58%  '
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = "'35.58"
$ws.Range('E14').Value = '  +0.92%  '
$ws.Range('D15').Value = '3.540.53'
$ws.Range('E15').Value = '  -1.80%  '
$ws.Range('D16').Value = '63.702.54'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '3.040.38'
$ws.Range('E17').Value = '  -2.01%  '
$ws.Range('E18').Value = '  +0.38%  '
$ws.Range('D19').Value = "'6.78"
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = "'474.39"
$ws.Range('E20').Value = '  -1.82%  '
$ws.Range('D21').Value = "'14.04"
$ws.Range('E21').Value = '  +1.58%  '
$ws.Range('D22').Value = "'14.67"
$ws.Range('E22').Value = '  +10.54%  '
$ws.Range('D23').Value = "'0.681"
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  -2.30%  '
$ws.Range('D25').Value = "'82.76"
$ws.Range('E25').Value = '  +2.37%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('E27').Value = '  -0.53%  '
$ws.Range('D28').Value = "'8.12"
$ws.Range('E28').Value = '  +0.52%  '
$ws.Range('E29').Value = '  -1.89%  '
$ws.Range('D30').Value = "'1.00"
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').Value = "'26.17"
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('E32').Value = '  -0.57%  '
$ws.Range('E33').Value = '  -1.15%  '
$ws.Range('D34').Value = "'5.75"
$ws.Range('E34').Value = '  +0.34%  '
$ws.Range('D35').Value = "'6.20"
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('D36').Value = "'54.64"
$ws.Range('E36').Value = '  -0.64%  '
$ws.Range('D37').Value = "'0.0408"
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').Value = "'440.23"
$ws.Range('E38').Value = '  -5.11%  '
$ws.Range('D39').Value = "'0.0811"
$ws.Range('E39').Value = '  -2.01%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '3.009.65'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = "'2.78"
$ws.Range('E41').Value = '  +2.34%  '
$ws.Range('E42').Value = '  +0.56%  '
$ws.Range('D43').Value = "'8.26"
$ws.Range('E43').Value = '  -0.18%  '
$ws.Range('D44').Value = "'0.269"
$ws.Range('E44').Value = '  +3.16%  '
$ws.Range('D45').Value = "'27.78"
$ws.Range('E45').Value = '  -2.44%  '
$ws.Range('D46').Value = "'2.25"
$ws.Range('E46').Value = '  +7.46%  '
$ws.Range('E48').Value = '  +0.76%  '
$ws.Range('D49').Value = "'118.41"
$ws.Range('E49').Value = '  -0.03%  '
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('D51').Value = "'2.08"
$ws.Range('E51').Value = '  +0.50%  '
